$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and D hold text that looks numeric/date-like ("2025-02-25",
# "08"); force a text number-format before assigning so Excel doesn't
# silently convert them to a date serial / drop the leading zero, then
# drop back to the Normal style so the new row keeps the same (unstyled)
# look as every other data row in the sheet.
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Value = "2025-02-25"
$ws.Range("A89").Style = "Normal"

$ws.Range("B89").Value = "23:01:58"
$ws.Range("C89").Value = "Tuesday"

$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "08"
$ws.Range("D89").Style = "Normal"

$ws.Range("E89").Value = 130729
$ws.Range("F89").Value = 141918
$ws.Range("G89").Value = 172790
$ws.Range("H89").Value = 159612
$ws.Range("I89").Value = -1
$ws.Range("J89").Value = 146655
$ws.Range("K89").Value = -1
$ws.Range("L89").Value = -1
$ws.Range("M89").Value = 193803
$ws.Range("N89").Value = 115384
$ws.Range("O89").Value = 46662
$ws.Range("P89").Value = 29475
$ws.Range("Q89").Value = 69258
$ws.Range("R89").Value = -1
$ws.Range("S89").Value = 49461
$ws.Range("T89").Value = -1
